$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restore the NFLX-only open/close/high/low/shares_outstanding/ticker values
# that were corrupted by data leaking in from other companies.
$rows = @(
    @{Row=2; D=59.64285659790039; E=79.5; F=82.30428314208984; G=58.46428680419922; H=423732334},
    @{Row=3; D=94.80571746826172; E=114.3099975585938; F=117.879997253418; G=92.28428649902344; H=423732334},
    @{Row=4; D=102.9100036621094; E=108.379997253418; F=115.8300018310547; G=96.26000213623048; H=423732334},
    @{Row=5; D=109; E=91.83999633789062; F=122.1800003051758; G=90.11000061035156; H=423732334},
    @{Row=6; D=102.9300003051758; E=90.02999877929688; F=111.8499984741211; G=88.20999908447266; H=423732334},
    @{Row=7; D=95; E=91.25; F=101.2699966430664; G=84.5; H=423732334},
    @{Row=8; D=98; E=124.870002746582; F=129.2899932861328; G=97.62999725341795; H=423732334},
    @{Row=9; D=124.9599990844727; E=140.7100067138672; F=143.4600067138672; G=124.3099975585938; H=423732334},
    @{Row=10; D=146.6999969482422; E=152.1999969482422; F=153.5200042724609; G=138.6600036621094; H=423732334},
    @{Row=11; D=149.8000030517578; E=181.6600036621093; F=191.5; G=144.25; H=423732334},
    @{Row=12; D=182.1100006103516; E=196.4299926757812; F=204.3800048828125; G=176.5800018310547; H=423732334},
    @{Row=13; D=196.1000061035156; E=270.2999877929688; F=286.8099975585937; G=195.4199981689453; H=423732334},
    @{Row=14; D=291.9400024414062; E=312.4599914550781; F=338.8200073242188; G=271.2200012207031; H=423732334},
    @{Row=15; D=385.4500122070313; E=337.4500122070312; F=419.7699890136719; G=328; H=423732334},
    @{Row=16; D=375.8500061035156; E=301.7799987792969; F=386.7999877929688; G=271.2099914550781; H=423732334},
    @{Row=17; D=259.2799987792969; E=339.5; F=358.8500061035156; G=256.5799865722656; H=423732334},
    @{Row=18; D=359; E=370.5400085449219; F=384.7999877929688; G=342.2699890136719; H=423732334},
    @{Row=19; D=373.5; E=322.989990234375; F=384.760009765625; G=305.8099975585937; H=423732334},
    @{Row=20; D=267.3500061035156; E=287.4100036621094; F=308.75; G=257.010009765625; H=423732334},
    @{Row=21; D=326.1000061035156; E=345.0899963378906; F=359.8500061035156; G=321.2000122070312; H=423732334},
    @{Row=22; D=376.0499877929688; E=419.8500061035156; F=449.5199890136719; G=357.510009765625; H=423732334},
    @{Row=23; D=454; E=488.8800048828125; F=575.3699951171875; G=454; H=423732334},
    @{Row=24; D=506.0299987792969; E=475.739990234375; F=572.489990234375; G=472.2099914550781; H=423732334},
    @{Row=25; D=539; E=532.3900146484375; F=593.2899780273438; G=485.6700134277344; H=423732334},
    @{Row=26; D=529.9299926757812; E=513.469970703125; F=563.5599975585938; G=499; H=423732334},
    @{Row=27; D=525.719970703125; E=517.5700073242188; F=557.5399780273438; G=504.6600036621094; H=423732334},
    @{Row=28; D=604.239990234375; E=690.3099975585938; F=690.969970703125; G=594.6799926757812; H=423732334},
    @{Row=29; D=605.6099853515625; E=427.1400146484375; F=609.989990234375; G=351.4599914550781; H=423732334},
    @{Row=30; D=376.7999877929688; E=190.3600006103516; F=396.0199890136719; G=185.6000061035156; H=423732334},
    @{Row=31; D=176.4900054931641; E=224.8999938964844; F=230.75; G=169.6999969482422; H=423732334},
    @{Row=32; D=237.8500061035156; E=291.8800048828125; F=305.6300048828125; G=211.729995727539; H=423732334},
    @{Row=33; D=298.0599975585937; E=353.8599853515625; F=369.0199890136719; G=288.7000122070312; H=423732334},
    @{Row=34; D=341.8299865722656; E=329.9299926757812; F=349.7999877929688; G=316.1000061035156; H=423732334},
    @{Row=35; D=439.760009765625; E=438.9700012207031; F=485; G=411.8800048828125; H=423732334},
    @{Row=36; D=377.4800109863281; E=411.6900024414063; F=418.8399963378906; G=344.7300109863281; H=423732334},
    @{Row=37; D=483.1900024414063; E=564.1099853515625; F=579.6400146484375; G=461.8599853515625; H=423732334},
    @{Row=38; D=608; E=550.6400146484375; F=639; G=542.010009765625; H=423732334},
    @{Row=39; D=674.8900146484375; E=628.3499755859375; F=697.489990234375; G=617; H=423732334},
    @{Row=40; D=713.6400146484375; E=756.030029296875; F=773; G=677.8800048828125; H=423732334},
    @{Row=41; D=895.5; E=976.760009765625; F=999; G=823.52001953125; H=423732334},
    @{Row=42; D=927.5; E=1131.719970703125; F=1133.199951171875; G=821.0999755859375; H=423732334},
    @{Row=43; D=1338.22998046875; E=1159.400024414062; F=1338.77001953125; G=1157.739990234375; H=423732334},
    @{Row=44; D=1175.489990234375; E=1094.68994140625; F=1248.599975585938; G=1094.410034179688; H=423732334}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = "NFLX"
}